$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: Activité text changes (value swap between 6/7 rows)
$ws.Range("B6").Value = "Modification des règles du jeu sur le repo distant"

# Row 7: Activité text + hours change
$ws.Range("B7").Value = "Test simple de l'application et de ses fonctionnalités"
$ws.Range("C7").Value = 3

# Row 8 (new data)
$ws.Range("A8").Value = 43004
$ws.Range("B8").Value = "Analyse de l'application : partie simple"
$ws.Range("C8").Value = 3

# Row 9 (new data)
$ws.Range("A9").Value = 43008
$ws.Range("B9").Value = "Analyse de l'application : partie tournoi"
$ws.Range("C9").Value = 2

# Row 11 (new data) - set before row 10 so shared-string order matches target
$ws.Range("A11").Value = 43016
$ws.Range("B11").Value = "Diagramme de classe côté serveur"
$ws.Range("C11").Value = 1.5

# Row 10 (new data)
$ws.Range("A10").Value = 43011
$ws.Range("B10").Value = "Début diagramme de classe (serveur, client et common)"
$ws.Range("C10").Value = 4

# Update the active selection to B9
$ws.Range("B9").Select()
